$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Footer "date" field: 6/22/2022 -> 6/24/2022 on the slide master and on
#    every slide layout (16 layouts + 1 master).
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "6/22/2022") {
                $sh.TextFrame.TextRange.Text = "6/24/2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShape $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 13 ("Conclusion") - two wording tweaks in the Content Placeholder.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange

$districtsPara = $tr13.Paragraphs(2, 1)
$districtsPara.Text = "Community Districts 111, 205, 206, 305, & 316 are the community districts that have the best options."

$demoPara = $tr13.Paragraphs(5, 1)
$demoPara.Text = "When it comes to demographics, the top community districts had more black, hispanic, and foreign born residents compared to the bottom districts which were more diverse. As previously stated, these results are most likely due to the city government having success getting underserved communities access to better food with the Shop Healthy NYC program. "

# ---------------------------------------------------------------------------
# 3. Slide 14 ("Things recommended for further investigation") - trim several
#    bullet points, resize/reposition the content placeholder, and drop the
#    extra line-spacing reduction on the autofit.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

# Delete from the highest paragraph index down so earlier indices stay valid.
$toDelete = @(13, 11, 10, 9, 8, 7, 3, 2)
foreach ($idx in $toDelete) {
    $tr14.Paragraphs($idx, 1).Delete()
}

$sh14.Left = 25.6696062992126
$sh14.Top = 75.75661417322834
$sh14.Width = 676.9029921259843
$sh14.TextFrame.AutoSize = 2
